$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating-point precision tweak on A11 (same timestamp, 18:00:16)
$ws.Range("A11").Value = 45863.75018738426

# Append the new row 12 captured by the scheduled task run
$ws.Range("A12").Value = 45863.79191518768
$ws.Range("B12").Value = 2025
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 15.56
$ws.Range("E12").Value = 81.63
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 4.87
$ws.Range("H12").Value = "E"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "19:00:21"

# Carry the date/time number format from the row above onto the new row's date cell
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat
